$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Smooth 6v"
$ws.Range("C11").Value = "17 Απρ"
$ws.Range("D11").Value = 2.8
$ws.Range("E11").Value = 9.07
$ws.Range("F11").Value = 13.84
$ws.Range("G11").Value = 'y="-math.sin(6*t-7.854)*20000-19980" ,t>0'

$ws.Range("B11:G11").HorizontalAlignment = -4108
$ws.Range("B11:G11").VerticalAlignment = -4108

$ws.Range("G11").Select()
